# update hotel reviews data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hotel_info")

# English_Reviews_num (G2) and Local_Rank (H2) for the hotel in row 2.
# These are text/shared-string cells in the source data (not numbers), so
# write them as formulas that yield text, then paste-special as values -
# this keeps the result stored as a string (t="s") without mutating the
# cell's number format / style, matching the source workbook exactly.
$ws.Range("G2").Formula = '="1"'
$ws.Range("H2").Formula = '="354"'
$ws.Range("G2:H2").Copy()
$ws.Range("G2:H2").PasteSpecial(-4163)
